$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 80
$ws.Range("I2").Value = 206
$ws.Range("J2").Value = 849
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 219
$ws.Range("N2").Value = 136
$ws.Range("P2").Value = 6
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 7
$ws.Range("S2").Value = 93
$ws.Range("T2").Value = 176
$ws.Range("U2").Value = 15
$ws.Range("V2").Value = 1338
$ws.Range("X2").Value = 1370
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 5
